$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(30, 1).Value = 29
$ws.Cells.Item(30, 2).Value = 2
$ws.Cells.Item(30, 3).Value = 31
$ws.Cells.Item(30, 4).Value = 10
$ws.Cells.Item(30, 5).Value = "System"
$ws.Cells.Item(30, 6).Value = "2025-03-03 17:38:18"
$ws.Cells.Item(30, 7).Value = 0

$ws.Cells.Item(31, 1).Value = 30
$ws.Cells.Item(31, 2).Value = 2
$ws.Cells.Item(31, 3).Value = 32
$ws.Cells.Item(31, 4).Value = 6
$ws.Cells.Item(31, 5).Value = "System"
$ws.Cells.Item(31, 6).Value = "2025-03-03 17:44:09"
$ws.Cells.Item(31, 7).Value = 0

$ws.Cells.Item(32, 1).Value = 31
$ws.Cells.Item(32, 2).Value = 2
$ws.Cells.Item(32, 3).Value = 33
$ws.Cells.Item(32, 4).Value = 20
$ws.Cells.Item(32, 5).Value = "System"
$ws.Cells.Item(32, 6).Value = "2025-03-03 17:45:14"
$ws.Cells.Item(32, 7).Value = 0

$ws.Cells.Item(33, 1).Value = 32
$ws.Cells.Item(33, 2).Value = 2
$ws.Cells.Item(33, 3).Value = 34
$ws.Cells.Item(33, 4).Value = 20
$ws.Cells.Item(33, 5).Value = "System"
$ws.Cells.Item(33, 6).Value = "2025-03-03 17:45:56"
$ws.Cells.Item(33, 7).Value = 0

$ws.Cells.Item(34, 1).Value = 33
$ws.Cells.Item(34, 2).Value = 2
$ws.Cells.Item(34, 3).Value = 35
$ws.Cells.Item(34, 4).Value = 4
$ws.Cells.Item(34, 5).Value = "System"
$ws.Cells.Item(34, 6).Value = "2025-03-03 17:49:17"
$ws.Cells.Item(34, 7).Value = 0

$ws.Cells.Item(35, 1).Value = 34
$ws.Cells.Item(35, 2).Value = 2
$ws.Cells.Item(35, 3).Value = 36
$ws.Cells.Item(35, 4).Value = 3
$ws.Cells.Item(35, 5).Value = "System"
$ws.Cells.Item(35, 6).Value = "2025-03-03 17:50:02"
$ws.Cells.Item(35, 7).Value = 0

$ws.Cells.Item(36, 1).Value = 35
$ws.Cells.Item(36, 2).Value = 2
$ws.Cells.Item(36, 3).Value = 37
$ws.Cells.Item(36, 4).Value = 3
$ws.Cells.Item(36, 5).Value = "System"
$ws.Cells.Item(36, 6).Value = "2025-03-03 17:53:52"
$ws.Cells.Item(36, 7).Value = 0

$ws.Cells.Item(37, 1).Value = 36
$ws.Cells.Item(37, 2).Value = 2
$ws.Cells.Item(37, 3).Value = 38
$ws.Cells.Item(37, 4).Value = 2
$ws.Cells.Item(37, 5).Value = "System"
$ws.Cells.Item(37, 6).Value = "2025-03-03 17:55:12"
$ws.Cells.Item(37, 7).Value = 0

$ws.Cells.Item(38, 1).Value = 37
$ws.Cells.Item(38, 2).Value = 2
$ws.Cells.Item(38, 3).Value = 39
$ws.Cells.Item(38, 4).Value = 3
$ws.Cells.Item(38, 5).Value = "System"
$ws.Cells.Item(38, 6).Value = "2025-03-03 18:01:42"
$ws.Cells.Item(38, 7).Value = 0

$ws.Cells.Item(39, 1).Value = 38
$ws.Cells.Item(39, 2).Value = 2
$ws.Cells.Item(39, 3).Value = 40
$ws.Cells.Item(39, 4).Value = 3
$ws.Cells.Item(39, 5).Value = "System"
$ws.Cells.Item(39, 6).Value = "2025-03-03 18:05:52"
$ws.Cells.Item(39, 7).Value = 0

$ws.Cells.Item(40, 1).Value = 39
$ws.Cells.Item(40, 2).Value = 2
$ws.Cells.Item(40, 3).Value = 41
$ws.Cells.Item(40, 4).Value = 1
$ws.Cells.Item(40, 5).Value = "System"
$ws.Cells.Item(40, 6).Value = "2025-03-03 18:12:41"
$ws.Cells.Item(40, 7).Value = 0

$ws.Cells.Item(41, 1).Value = 40
$ws.Cells.Item(41, 2).Value = 8
$ws.Cells.Item(41, 3).Value = 42
$ws.Cells.Item(41, 4).Value = 1
$ws.Cells.Item(41, 5).Value = "System"
$ws.Cells.Item(41, 6).Value = "2025-03-03 18:19:56"
$ws.Cells.Item(41, 7).Value = 0
